$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 57, shifting existing rows 57:67 down to 58:68.
$ws.Rows(57).Insert()

# Populate the newly inserted row 57 with the new record's data.
$ws.Range("A57").Value = 5
$ws.Range("B57").Value = "Macroferia Regional de Talca"
$ws.Range("C57").Value = "Maule"
$ws.Range("D57").Value = 44637
$ws.Range("E57").Value = 7
$ws.Range("F57").Value = "Fruta"
$ws.Range("G57").Value = 100101
$ws.Range("H57").Value = "Berries"
$ws.Range("I57").Value = 100101001
$ws.Range("J57").Value = "Arándano (blue)"
$ws.Range("K57").Value = "Sin especificar"
$ws.Range("L57").Value = "Segunda"
$ws.Range("M57").Value = 150
$ws.Range("N57").Value = 3000
$ws.Range("O57").Value = 3000
$ws.Range("P57").Value = 3000
$ws.Range("Q57").Value = '$/bandeja 2 kilos'
$ws.Range("R57").Value = "Provincia de Linares"
$ws.Range("S57").Value = 1500
$ws.Range("T57").Value = 2
